$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new Price text, new Volume(1h) text) updates scraped on 2024-01-05
$updates = @(
    @{ Row = 2; D = "43.977.03"; E = "  +1.23%  " }
    @{ Row = 3; D = "2.241.34"; E = "  -0.01%  " }
    @{ Row = 4; D = $null; E = "  +0.20%  " }
    @{ Row = 5; D = "317.18"; E = "  -0.36%  " }
    @{ Row = 6; D = "100.34"; E = "  +0.09%  " }
    @{ Row = 7; D = "0.571"; E = "  -2.08%  " }
    @{ Row = 8; D = $null; E = "  +0.16%  " }
    @{ Row = 9; D = "0.540"; E = "  -4.33%  " }
    @{ Row = 10; D = "36.59"; E = "  -1.78%  " }
    @{ Row = 11; D = "0.0824"; E = "  -1.21%  " }
    @{ Row = 12; D = "7.48"; E = "  -3.28%  " }
    @{ Row = 13; D = $null; E = "  -2.19%  " }
    @{ Row = 14; D = "2.583.43"; E = "  +0.11%  " }
    @{ Row = 15; D = "0.843"; E = "  -2.94%  " }
    @{ Row = 16; D = "2.242.73"; E = "  +0.59%  " }
    @{ Row = 17; D = "14.09"; E = "  -2.01%  " }
    @{ Row = 18; D = "43.874.06"; E = "  +1.12%  " }
    @{ Row = 19; D = "13.15"; E = "  -7.78%  " }
    @{ Row = 20; D = "0.0₃0967"; E = "  -0.72%  " }
    @{ Row = 21; D = "6.40"; E = "  -4.11%  " }
    @{ Row = 22; D = "65.25"; E = "  -0.37%  " }
    @{ Row = 23; D = "3.07"; E = "  -4.13%  " }
    @{ Row = 24; D = "234.30"; E = "  -1.15%  " }
    @{ Row = 25; D = "2.04"; E = "  -6.20%  " }
    @{ Row = 26; D = $null; E = "  +0.15%  " }
    @{ Row = 27; D = "10.39"; E = "  +3.01%  " }
    @{ Row = 28; D = $null; E = "  -0.52%  " }
    @{ Row = 29; D = "37.10"; E = "  +1.09%  " }
    @{ Row = 30; D = "6.10"; E = "  -4.92%  " }
    @{ Row = 31; D = "158.83"; E = "  +0.44%  " }
    @{ Row = 32; D = "19.98"; E = "  -1.71%  " }
    @{ Row = 33; D = "0.0842"; E = "  -3.45%  " }
    @{ Row = 34; D = $null; E = "  -1.66%  " }
    @{ Row = 35; D = "3.16"; E = "  -1.01%  " }
    @{ Row = 36; D = "0.112"; E = "  +7.34%  " }
    @{ Row = 37; D = "1.92"; E = "  +1.14%  " }
    @{ Row = 38; D = $null; E = "  -2.51%  " }
    @{ Row = 39; D = "16.00"; E = "  +10.05%  " }
    @{ Row = 40; D = "3.63"; E = "  -2.57%  " }
    @{ Row = 41; D = "4.08"; E = "  -7.16%  " }
    @{ Row = 42; D = $null; E = "  -3.41%  " }
    @{ Row = 43; D = $null; E = "  +0.17%  " }
    @{ Row = 44; D = "1.741.28"; E = "  -5.27%  " }
    @{ Row = 45; D = "0.196"; E = "  -3.60%  " }
    @{ Row = 46; D = "81.30"; E = "  -4.02%  " }
    @{ Row = 47; D = "74.03"; E = "  -0.57%  " }
    @{ Row = 48; D = "5.13"; E = "  -3.63%  " }
    @{ Row = 49; D = "101.68"; E = "  -1.97%  " }
    @{ Row = 50; D = $null; E = "  +1.68%  " }
    @{ Row = 51; D = "57.02"; E = "  -2.46%  " }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # Force text so Excel does not reinterpret values like "317.18" as numbers,
        # matching the inline-string storage used by the source sheet, then restore
        # the cell style so no stray formatting is introduced.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
